$d = $word.ActiveDocument

# Helper: split the run that contains character position $pos into two runs by
# inserting and then immediately removing a bookmark there. Word (and this
# COM-interop runtime) breaks the run at the insertion point when a bookmark
# is added mid-run; deleting the bookmark again leaves the split runs intact.
function SplitRunAt($pos) {
    $pt = $d.Range($pos, $pos)
    $d.Bookmarks.Add("__TempSplitMark__", $pt)
    $d.Bookmarks("__TempSplitMark__").Delete()
}

# ---------------------------------------------------------------------------
# Change 1: remove the stray "_GoBack" bookmark that sits right after
# "The force will be with you, always."
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# Change 2: the "_GoBack" bookmark now belongs in the middle of the sentence
# about passing a random number in as the index - right after "...random num"
# and before "ber for the index...". Re-add it there (this naturally splits
# the run in two, matching the target XML).
# ---------------------------------------------------------------------------
$rNum = $d.Content
$rNum.Find.Execute(", passing in a random num", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$goBackPoint = $d.Range($rNum.End, $rNum.End)
$d.Bookmarks.Add("_GoBack", $goBackPoint)

# ---------------------------------------------------------------------------
# Change 3: in the "Luke Skywalker Level" stretch paragraph, rename the
# referenced button from "Force Sight" to "Force Read" (three occurrences),
# leaving each occurrence of the word isolated in its own run, just like the
# existing "Force Read" button description elsewhere in the document.
# ---------------------------------------------------------------------------

# Locate the paragraph that still mentions "Force Sight".
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Force Sight*") {
        $targetPara = $p
    }
}

if ($targetPara -ne $null) {
    $paraStart = $targetPara.Range.Start

    # Replace the word everywhere in this paragraph (Word merges the runs of
    # the edited paragraph in the process - that's fine, we restore the run
    # boundaries explicitly afterwards).
    $rParaText = $d.Range($paraStart, $targetPara.Range.End)
    $rParaText.Find.Execute("Sight", $true, $false, $false, $false, $false, $true, 1, $false, "Read", 2) | Out-Null

    # Re-establish the run split right before "Change ..." (the sentence
    # preceding it, "You will need to complete Obi-wan Kenobi level first. ",
    # stays in its own run).
    $rBoundary = $d.Content
    $rBoundary.Start = $paraStart
    $rBoundary.Find.Execute("You will need to complete Obi-wan Kenobi level first. ", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
    $changeBoundary = $rBoundary.End

    $splitPositions = New-Object System.Collections.ArrayList
    [void]$splitPositions.Add($changeBoundary)

    # Find each "Read" occurrence (there are three) and remember the
    # positions right before and right after it so it ends up isolated in
    # its own run.
    $searchStart = $changeBoundary
    for ($i = 0; $i -lt 3; $i++) {
        $rWord = $d.Content
        $rWord.Start = $searchStart
        $rWord.Find.Execute("Read", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
        [void]$splitPositions.Add($rWord.Start)
        [void]$splitPositions.Add($rWord.End)
        $searchStart = $rWord.End
    }

    foreach ($p in $splitPositions) {
        SplitRunAt($p)
    }
}
